$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A24").Value = 43898
$ws.Range("B24").Value = "老王 "
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = "Data preparing/Final feature engineering_v0.1/Model selection_v0.1"
$ws.Range("E24").Value = "Please Review the comment part on Feature engineering part and provide feedback"

$ws.Rows.Item(24).RowHeight = 28.8

$wb.Save()
